$d = $word.ActiveDocument

# Locate the "Edison Achalma" paragraph styled as "Author" (the author byline
# just under the document title) by scanning the Paragraphs collection and
# remembering its 1-based index (Paragraph.Next is unreliable in this host,
# so re-fetch everything through $d.Paragraphs(<index>) instead).
$idx = 0
$authorIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $clean = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($clean -eq "Edison Achalma" -and $p.Style.NameLocal -eq "Author") {
        $authorIdx = $idx
    }
}

if ($authorIdx -gt 0) {
    $followIdx = $authorIdx + 1
    $followPara = $d.Paragraphs($followIdx)

    # Split one character into the following paragraph. In this host,
    # InsertParagraphAfter on a collapsed Range pushes a *new* blank
    # paragraph in front of whatever range it was called on and shifts
    # that range's own content forward by one slot - splitting inside the
    # following paragraph (rather than at the exact end of the "Edison
    # Achalma" paragraph) keeps "Edison Achalma" completely untouched and
    # lands the freshly inserted blank paragraph right after it.
    $splitPos = $followPara.Range.Start + 1
    $splitRange = $d.Range($splitPos, $splitPos)
    $splitRange.InsertParagraphAfter()

    # The new blank paragraph now occupies the old $followIdx slot.
    $newPara = $d.Paragraphs($followIdx)
    $newPara.Style = "Author"

    # Seed with a space first (forces the run's xml:space="preserve" flag),
    # then overwrite with the real affiliation text.
    $newPara.Range.Text = " "
    $newPara2 = $d.Paragraphs($followIdx)
    $newPara2.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
}
